# Update "想去人数" (want-to-go count) values in the F column for both the
# "展览" sheet and the "全部类型" sheet, reflecting newly scraped numbers.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 228
$ws1.Range("F3").Value = 1073
$ws1.Range("F8").Value = 49
$ws1.Range("F9").Value = 6679
$ws1.Range("F15").Value = 1068
$ws1.Range("F16").Value = 15980
$ws1.Range("F22").Value = 11238
$ws1.Range("F24").Value = 4416

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 228
$ws4.Range("F3").Value = 1073
$ws4.Range("F9").Value = 49
$ws4.Range("F10").Value = 6679
$ws4.Range("F17").Value = 1068
$ws4.Range("F18").Value = 15980
$ws4.Range("F25").Value = 11238
$ws4.Range("F27").Value = 4416
